$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "purpose" column (E) for all sample rows changes from "S.GISH" to "fullRNASEQ"
$ws.Range("E2:E24").Value = "fullRNASEQ"

# Leave the selection on the last edited cell, matching the saved view state
[void]$ws.Range("E24").Select()
